$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (values shuffled among B/C/D)
$ws.Range("B2").Value = 12579507344534390
$ws.Range("C2").Value = 12579507344534380
$ws.Range("D2").Value = 12579507344534380

# Row 3 - RandomForestRegressor (values replaced)
$ws.Range("B3").Value = 2234238614511.596
$ws.Range("C3").Value = 2217121683463.108
$ws.Range("D3").Value = 820801723566729.6

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 2302415974387.048
$ws.Range("C4").Value = 2159218901556.084
$ws.Range("D4").Value = 247277787747725.9

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 333746347753674.9
$ws.Range("C5").Value = 1787343908159784
$ws.Range("D5").Value = 5355413201749120
